$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @{ row=259; B="Al25Hf25Nb25Ti25"; J=1700000000 },
  @{ row=260; B="Al20Hf24Nb29Ti27"; J=1600000000 }
)

foreach ($item in $data) {
  $ws.Cells.Item($item.row, 2).Value = $item.B
  $ws.Cells.Item($item.row, 10).Value = $item.J
}
Write-Output ($ws.Cells.Item(259,2).Value())
Write-Output ($ws.Cells.Item(259,10).Value())
Write-Output ($ws.Cells.Item(260,2).Value())
Write-Output ($ws.Cells.Item(260,10).Value())
